$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 58: fill in the new journal entry (date, description, hours)
$ws.Range("A58").Value = 43227
$ws.Range("B58").Value = "Préparation des méthodes pour transaction, résolutions de certains problèmes avec l'équipe et Guillaume."
$ws.Range("C58").Value = 1.5

# Row 58 grew a bit taller to fit the wrapped description text
$ws.Rows(58).RowHeight = 30

# Row 59: just the date carried over, rest left blank
$ws.Range("A59").Value = 43227

# Move the active selection to B59, matching where the user ended up
[void]$ws.Range("B59").Select()

Write-Output "done"
